# Updated cryptos list (Price + Volume(1h) columns) for rows 2-51.
# D-column price cells are plain text (e.g. "26.876.27", "1.000") so we
# force NumberFormat "@" before assigning to stop Excel from coercing
# them into numbers (which would strip significant trailing zeros or
# mis-parse multi-dot values). E-column volume strings already contain
# non-numeric padding/spacing so they stay text automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.876.27'
$ws.Range("E2").Value = '  -0.85%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.814.59'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.65'
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4679'
$ws.Range("E7").Value = '  +1.19%  '
$ws.Range("E8").Value = '  -1.54%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07372'
$ws.Range("E9").Value = '  -0.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8703'
$ws.Range("E10").Value = '  +0.92%  '
$ws.Range("E11").Value = '  -0.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.836.91'
$ws.Range("E12").Value = '  +1.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.364'
$ws.Range("E13").Value = '  -0.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.34'
$ws.Range("E14").Value = '  +0.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07083'
$ws.Range("E15").Value = '  +0.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.500'
$ws.Range("E16").Value = '  -2.27%  '
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008709'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("E20").Value = '  -0.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.894.76'
$ws.Range("E21").Value = '  -0.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.333'
$ws.Range("E22").Value = '  +0.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.53'
$ws.Range("E23").Value = '  -2.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.039.98'
$ws.Range("E24").Value = '  -0.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.893'
$ws.Range("E25").Value = '  -1.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.77'
$ws.Range("E26").Value = '  +0.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.189'
$ws.Range("E27").Value = '  -0.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.35'
$ws.Range("E28").Value = '  -0.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.296'
$ws.Range("E29").Value = '  +0.63%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.51'
$ws.Range("E30").Value = '  -1.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08921'
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7648'
$ws.Range("E32").Value = '  -0.98%  '
$ws.Range("E33").Value = '  -0.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.480'
$ws.Range("E34").Value = '  -1.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.920'
$ws.Range("E35").Value = '  +0.75%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.001'
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.093'
$ws.Range("E37").Value = '  -3.15%  '
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05269'
$ws.Range("E39").Value = '  +0.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.943'
$ws.Range("E40").Value = '  +0.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.237'
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5335'
$ws.Range("E42").Value = '  +1.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.358'
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1661'
$ws.Range("E44").Value = '  -0.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.444'
$ws.Range("E45").Value = '  -1.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4936'
$ws.Range("E46").Value = '  -1.78%  '
$ws.Range("E47").Value = '  +1.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.670'
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.65'
$ws.Range("E50").Value = '  -2.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06281'
$ws.Range("E51").Value = '  -0.62%  '
